# Added validation for add bids: cannot bid for the same course in a single
# bid. Append the two new bug-log rows (S/N 21 and a follow-up entry) to the
# "Bug Log" sheet, matching the formatting of the surrounding rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bug Log")

# --- Row 23 ------------------------------------------------------------
# Same visual layout as row 22 (S/N 20), so clone its cell formatting first.
$ws.Range("A22:H22").Copy()
$ws.Range("A23:H23").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 24 --------------------------------------------------------------
# Only columns C:H are populated (no S/N or Iteration), matching the layout
# used by row 19. Column H there is empty, so borrow that single cell's
# format from row 20 (same "blank result" style used elsewhere).
$ws.Range("C19:G19").Copy()
$ws.Range("C24:G24").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("H20").Copy()
$ws.Range("H24").PasteSpecial(-4122)       # xlPasteFormats

$excel.CutCopyMode = 0

# --- Values for row 23 (S/N 21) ------------------------------------------
$ws.Cells.Item(23, 1).Value2 = 21
$ws.Cells.Item(23, 2).Value2 = 3
$ws.Cells.Item(23, 3).Value2 = "Add bid "
$ws.Cells.Item(23, 4).Value2 = "No validation if you bid for two sections of the same course at the same time"
$ws.Cells.Item(23, 5).Value2 = "Resolved"
$ws.Cells.Item(23, 6).Value2 = 43778
$ws.Cells.Item(23, 7).Value2 = 43778
$ws.Cells.Item(23, 8).Value2 = "Matthew & Sheng Qin"

# --- Values for row 24 (follow-up bug, still unresolved) ------------------
$ws.Cells.Item(24, 3).Value2 = "Add bid "
$ws.Cells.Item(24, 4).Value2 = "No validation if you bid for two sections of the same course in separate bids"
$ws.Cells.Item(24, 5).Value2 = "Unresolved"
$ws.Cells.Item(24, 6).Value2 = 43778

# --- Row heights matching the source rows they were cloned from -----------
$ws.Rows.Item(23).RowHeight = 26.25
$ws.Rows.Item(24).RowHeight = 15.75

# --- Update the view state (scroll position / active selection) -----------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 16
$win.ScrollColumn = 3
$ws.Range("H24").Select()
